$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab Sheet1 -> tests
$ws.Name = "tests"

# Bold the header row (A1:B1)
$ws.Range("A1:B1").Font.Bold = $true

# Column A/B widen slightly to fit the new bold header text (values chosen
# to land as close as possible to the widths Excel's own "best fit"
# recalculation produced: ~13.44 chars for column A, ~9.22 chars for column B).
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(2).ColumnWidth = 8.333333333333334

# Move the active selection to G7
[void]$ws.Range("G7").Select()
